# Update column G ("K") values per regenerated save_data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 3
$ws.Range("G7").Value = 1
$ws.Range("G8").Value = 1
